$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds numeric-looking text values (e.g. " 75") that must stay text.
# Excel's Value setter auto-converts plain numeric strings to numbers, so we
# temporarily force a text number format, assign the value, then clear the
# format again so no residual style index is left on the cell.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = " 77"
$ws.Range("A2").ClearFormats()

$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = " 23"
$ws.Range("A3").ClearFormats()

# Row 2 numeric updates
$ws.Range("B2").Value = 222
$ws.Range("C2").Value = 1.01
$ws.Range("D2").Value = 203.1
$ws.Range("E2").Value = 1.04
$ws.Range("F2").Value = 140.2
$ws.Range("G2").Value = 0.72
$ws.Range("H2").Value = 0.73
$ws.Range("L2").Value = 162
$ws.Range("M2").Value = 156
$ws.Range("N2").Value = 73
$ws.Range("O2").Value = 13
$ws.Range("P2").Value = 22
